$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '43.534.52'
$ws.Range("E2").Value = '  +2.62%  '

$ws.Range("D3").Value = '2.315.51'
$ws.Range("E3").Value = '  +1.81%  '

$ws.Range("E4").Value = '  +0.06%  '

$ws.Range("D5").Value = '311.79'
$ws.Range("E5").Value = '  +1.60%  '

$ws.Range("D6").Value = '102.76'
$ws.Range("E6").Value = '  +5.32%  '

$ws.Range("E7").Value = '  +1.72%  '

$ws.Range("E8").Value = '  +0.01%  '

$ws.Range("D9").Value = '0.533'
$ws.Range("E9").Value = '  +7.59%  '

$ws.Range("D10").Value = '35.98'
$ws.Range("E10").Value = '  +1.53%  '

$ws.Range("D11").Value = '0.0818'
$ws.Range("E11").Value = '  +3.11%  '

$ws.Range("E12").Value = '  -0.46%  '

$ws.Range("D13").Value = '7.07'
$ws.Range("E13").Value = '  +2.42%  '

$ws.Range("D14").Value = '2.674.09'
$ws.Range("E14").Value = '  +1.88%  '

$ws.Range("D15").Value = '15.04'
$ws.Range("E15").Value = '  +1.32%  '

$ws.Range("D16").Value = '2.309.67'
$ws.Range("E16").Value = '  +0.59%  '

$ws.Range("D17").Value = '0.812'
$ws.Range("E17").Value = '  +1.99%  '

$ws.Range("D18").Value = '43.440.36'
$ws.Range("E18").Value = '  +2.85%  '

$ws.Range("D19").Value = '12.49'
$ws.Range("E19").Value = '  -0.21%  '

$ws.Range("D20").Value = '0.0₃0930'
$ws.Range("E20").Value = '  +2.46%  '

$ws.Range("E21").Value = '  +1.77%  '

$ws.Range("D22").Value = '68.32'
$ws.Range("E22").Value = '  +0.09%  '

$ws.Range("D23").Value = '242.37'
$ws.Range("E23").Value = '  +1.41%  '

$ws.Range("E24").Value = '  +4.71%  '

$ws.Range("E25").Value = '  +1.74%  '

$ws.Range("E26").Value = '  -0.10%  '

$ws.Range("D28").Value = '24.80'
$ws.Range("E28").Value = '  +4.67%  '

$ws.Range("B29").Value = 'Toncoin'
$ws.Range("C29").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D29").Value = '2.25'
$ws.Range("E29").Value = '  +6.38%  '

$ws.Range("B30").Value = 'InjectiveProtocol'
$ws.Range("C30").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D30").Value = '37.12'
$ws.Range("E30").Value = '  -2.54%  '

$ws.Range("B31").Value = 'Cosmos'
$ws.Range("C31").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D31").Value = '9.66'
$ws.Range("E31").Value = '  +1.27%  '

$ws.Range("D32").Value = '168.31'
$ws.Range("E32").Value = '  +3.92%  '

$ws.Range("D33").Value = '5.31'
$ws.Range("E33").Value = '  +0.95%  '

$ws.Range("E34").Value = '  +0.04%  '

$ws.Range("D35").Value = '3.12'
$ws.Range("E35").Value = '  -1.96%  '

$ws.Range("E36").Value = '  +6.79%  '

$ws.Range("D37").Value = '0.0747'
$ws.Range("E37").Value = '  +0.96%  '

$ws.Range("D38").Value = '17.58'
$ws.Range("E38").Value = '  +0.18%  '

$ws.Range("E39").Value = '  +1.80%  '

$ws.Range("D40").Value = '1.88'
$ws.Range("E40").Value = '  +2.56%  '

$ws.Range("D41").Value = '0.116'
$ws.Range("E41").Value = '  +1.80%  '

$ws.Range("D42").Value = '4.34'
$ws.Range("E42").Value = '  +5.97%  '

$ws.Range("E43").Value = '  -0.38%  '

$ws.Range("D44").Value = '19.44'
$ws.Range("E44").Value = '  +2.12%  '

$ws.Range("D45").Value = '0.0291'
$ws.Range("E45").Value = '  +3.11%  '

$ws.Range("D46").Value = '1.972.86'
$ws.Range("E46").Value = '  +1.04%  '

$ws.Range("D47").Value = '3.01'
$ws.Range("E47").Value = '  +3.52%  '

$ws.Range("E48").Value = '  -0.52%  '

$ws.Range("D49").Value = '55.66'
$ws.Range("E49").Value = '  +3.59%  '

$ws.Range("E50").Value = '  +6.41%  '

$ws.Range("D51").Value = '1.57'
$ws.Range("E51").Value = '  +6.57%  '
